$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.868.34"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "1.620.17"
$ws.Range("E3").Value = "  -2.12%  "

$ws.Range("E4").Value = "  -1.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.83"
$ws.Range("E5").Value = "  -2.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  -1.97%  "

$ws.Range("E7").Value = "  -0.95%  "

$ws.Range("E8").Value = "  -1.99%  "

$ws.Range("E9").Value = "  -3.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.48"
$ws.Range("E10").Value = "  -5.51%  "

$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("D12").Value = "1.845.15"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.624.20"
$ws.Range("E13").Value = "  -4.14%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -3.18%  "

$ws.Range("E15").Value = "  -3.68%  "

$ws.Range("D16").Value = "25.878.08"
$ws.Range("E16").Value = "  -0.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.51"
$ws.Range("E17").Value = "  -3.46%  "

$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("E18").Value = "  -3.60%  "

$ws.Range("E19").Value = "  -1.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.54"
$ws.Range("E20").Value = "  -1.70%  "

$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.50"
$ws.Range("E22").Value = "  -3.20%  "

$ws.Range("E23").Value = "  -2.91%  "

$ws.Range("E24").Value = "  +2.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.93"
$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("E26").Value = "  -1.25%  "

$ws.Range("E27").Value = "  -4.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.64"
$ws.Range("E28").Value = "  -3.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.20"
$ws.Range("E29").Value = "  -2.63%  "

$ws.Range("E30").Value = "  -1.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0478"
$ws.Range("E31").Value = "  -2.57%  "

$ws.Range("E32").Value = "  -4.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.10"
$ws.Range("E33").Value = "  -5.81%  "

$ws.Range("E34").Value = "  -2.96%  "

$ws.Range("E35").Value = "  -3.43%  "

$ws.Range("D36").Value = "1.125.46"
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.840"
$ws.Range("E37").Value = "  -7.31%  "

$ws.Range("E38").Value = "  -4.47%  "

$ws.Range("E39").Value = "  -2.78%  "

$ws.Range("E40").Value = "  -4.59%  "

$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").Value = "1.755.56"
$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.749"
$ws.Range("E43").Value = "  -6.70%  "

$ws.Range("E44").Value = "  -5.75%  "

$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.51"
$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.02"
$ws.Range("E47").Value = "  -4.58%  "

$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.411"
$ws.Range("E49").Value = "  -2.15%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.47"
$ws.Range("E50").Value = "  -3.57%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.90%  "

